$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 6 (shifts the old row 6 down to row 8),
# then insert one more row after that (new row 9).
$ws.Rows("6:7").Insert()
$ws.Rows("9:9").Insert()

# Row 6: new Mischfond product
$ws.Range("A6").Value = "PIB_VRWestmuensterland_MischfondNachhaltig_623669.pdf"
$ws.Range("B6").Value = "Mischfond"
$ws.Range("C6").Value = 623669
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = "mittelfristig"
$ws.Range("F6").Value = "mittleres Risiko"
$ws.Range("G6").Value = "ja"
$ws.Range("H6").Value = "ja"

# Row 7: new Privatfond product
$ws.Range("A7").Value = "PIB_Union_PrivatFond_12345678.pdf"
$ws.Range("B7").Value = "Privatfond"
$ws.Range("C7").Value = 12345678
$ws.Range("D7").Value = 10000
$ws.Range("E7").Value = "mittelfristig"
$ws.Range("F7").Value = "mittleres Risiko"
$ws.Range("G7").Value = "ja"
$ws.Range("H7").Value = "nein"

# Row 9: new Aktienfond product
$ws.Range("A9").Value = "PIB_UnionAsia_Aktienfond_971267.pdf"
$ws.Range("B9").Value = "Aktienfond"
$ws.Range("C9").Value = 971267
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = "langfristig"
$ws.Range("F9").Value = "hohes Risiko"
$ws.Range("G9").Value = "ja"
$ws.Range("H9").Value = "nein"

# Re-fit column widths now that longer text has been entered (matches the
# bestFit-recalculated widths for columns A, B and F in the target file).
# The host engine stores ColumnWidth assignments with a fixed +5/6 character
# padding added on save, so the input is pre-compensated by that amount to
# land as close as possible to the target stored width.
$padding = 5 / 6
$ws.Columns("A:A").ColumnWidth = 51.44140625 - $padding
$ws.Columns("B:B").ColumnWidth = 16.5546875 - $padding
$ws.Columns("F:F").ColumnWidth = 13.44140625 - $padding

# Final selection lands on A9, matching the author's last action.
$ws.Range("A9").Select() | Out-Null
